$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B column (DOF count) for specific rows per the diff
$ws.Cells.Item(15, 2).Value = 1
$ws.Cells.Item(21, 2).Value = 4
$ws.Cells.Item(22, 2).Value = 3
$ws.Cells.Item(23, 2).Value = 1
$ws.Cells.Item(24, 2).Value = 1
$ws.Cells.Item(27, 2).Value = 4
$ws.Cells.Item(28, 2).Value = 4
$ws.Cells.Item(30, 2).Value = 3
$ws.Cells.Item(34, 2).Value = 5
$ws.Cells.Item(35, 2).Value = 4
$ws.Cells.Item(36, 2).Value = 5

# Update H and I columns (all rows 1-36)
for ($r = 1; $r -le 36; $r++) {
    $ws.Cells.Item($r, 8).Value = 25.3
    $ws.Cells.Item($r, 9).Value = 25.3
}
